$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 16.084
$ws.Range("B12").Value = 5.2095
$ws.Range("E14").Value = 16.60170000000001
$ws.Range("E26").Value = 16.0452
$ws.Range("B27").Value = 5.583600000000001
$ws.Range("E31").Value = 16.2697
$ws.Range("B32").Value = 6.5223
$ws.Range("E35").Value = 16.7867
$ws.Range("B36").Value = 8.662600000000001
$ws.Range("E37").Value = 16.71190000000001
$ws.Range("B38").Value = 4.8307
$ws.Range("E45").Value = 16.5194
$ws.Range("B46").Value = 5.673700000000004
$ws.Range("E52").Value = 17.21020000000001
$ws.Range("B54").Value = 4.439799999999998
$ws.Range("B55").Value = 5.573999999999999
$ws.Range("B56").Value = 5.130599999999998
$ws.Range("E57").Value = 16.67680000000001
$ws.Range("B67").Value = 5.445899999999998
$ws.Range("B69").Value = 5.410599999999999
$ws.Range("B72").Value = 5.228000000000004
$ws.Range("E81").Value = 16.4836
$ws.Range("B83").Value = 5.742499999999998
$ws.Range("E83").Value = 16.51129999999999
$ws.Range("B86").Value = 4.959200000000003
$ws.Range("B91").Value = 5.486200000000002
$ws.Range("B93").Value = 5.921700000000004
$ws.Range("B99").Value = 4.747999999999998
$ws.Range("E100").Value = 16.4709
$ws.Range("E102").Value = 16.79179999999999

$wb.Save()
